$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Exported On" timestamp, preserving the bold "Exported On: " run
# and only replacing the date/time text of the second (non-bold) run.
$c = $ws.Range("A2").Characters(14, 20)
$c.Text = "06/17/2025 08:14 PM"

# Update Total Cost (F) and Quantity (K) values for the affected inventory rows.
$ws.Range("F11").Value = 8.48
$ws.Range("K11").Value = 8.0
$ws.Range("F14").Value = 10.33
$ws.Range("K14").Value = 1.0
$ws.Range("F15").Value = 16.29
$ws.Range("K15").Value = 9.0
$ws.Range("F18").Value = 0.0
$ws.Range("K18").Value = -1.0
$ws.Range("F22").Value = 0.0
$ws.Range("K22").Value = -1.0
$ws.Range("F25").Value = 8.71
$ws.Range("K25").Value = 13.0
$ws.Range("F28").Value = 9.0
$ws.Range("K28").Value = 15.0
$ws.Range("F32").Value = 320.7
$ws.Range("K32").Value = 10.0
$ws.Range("F34").Value = 35.52
$ws.Range("K34").Value = 32.0
$ws.Range("F39").Value = 24.48
$ws.Range("K39").Value = 12.0
$ws.Range("F43").Value = 106.98
$ws.Range("K43").Value = 2.0
$ws.Range("F44").Value = 68.46
$ws.Range("K44").Value = 2.0
$ws.Range("F45").Value = 482.35
$ws.Range("K45").Value = 11.0
$ws.Range("F49").Value = 211.59
$ws.Range("K49").Value = 9.0
$ws.Range("F50").Value = 12.45
$ws.Range("K50").Value = 5.0
$ws.Range("F54").Value = 7.68
$ws.Range("K54").Value = 4.0
$ws.Range("F55").Value = 20.88
$ws.Range("K55").Value = 9.0
$ws.Range("F62").Value = 70.0
$ws.Range("K62").Value = 2.0
$ws.Range("F99").Value = 97.23
$ws.Range("K99").Value = 7.0
$ws.Range("F114").Value = 262.15
$ws.Range("K114").Value = 7.0
$ws.Range("F118").Value = 272.85
$ws.Range("K118").Value = 1.0
$ws.Range("F127").Value = 47.04
$ws.Range("K127").Value = 4.0
$ws.Range("F130").Value = 254.66
$ws.Range("K130").Value = 2.0
$ws.Range("F131").Value = 509.32
$ws.Range("K131").Value = 4.0
$ws.Range("F133").Value = 60.9
$ws.Range("K133").Value = 6.0
$ws.Range("F134").Value = 42.72
$ws.Range("K134").Value = 4.0
$ws.Range("F136").Value = 20.3
$ws.Range("K136").Value = 2.0
$ws.Range("F140").Value = 38.46
$ws.Range("K140").Value = 3.0
$ws.Range("F147").Value = 15.9
$ws.Range("K147").Value = 5.0
$ws.Range("F151").Value = 25.65
$ws.Range("K151").Value = 3.0
$ws.Range("F159").Value = 42.78
$ws.Range("K159").Value = 2.0
$ws.Range("F160").Value = 38.46
$ws.Range("K160").Value = 6.0
$ws.Range("F164").Value = 11.79
$ws.Range("K164").Value = 9.0
$ws.Range("F168").Value = 9.28
$ws.Range("K168").Value = 4.0
$ws.Range("F170").Value = 77.0
$ws.Range("K170").Value = 4.0
$ws.Range("F171").Value = 2.56
$ws.Range("K171").Value = 2.0
$ws.Range("F173").Value = 4.9
$ws.Range("K173").Value = 2.0
$ws.Range("F183").Value = 50767.37
$ws.Range("K183").Value = 955.0
$ws.Range("F185").Value = 50767.37
$ws.Range("K185").Value = 955.0
